$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell reference -> new text value.
# Source data keeps these columns as plain text (Price / Volume% / Hora),
# so we force text via NumberFormat "@" before writing, then restore the
# original (General/style-0) number format by pasting formats only from
# column F of the same row (already plain text, unstyled).
$edits = @(
    @{Cell="D2"; Value="303.84"},
    @{Cell="E2"; Value="3.84%"},
    @{Cell="G2"; Value="18"},
    @{Cell="D3"; Value="35.76"},
    @{Cell="E3"; Value="15.21%"},
    @{Cell="G3"; Value="18"},
    @{Cell="D4"; Value="5.072"},
    @{Cell="E4"; Value="2.20%"},
    @{Cell="G4"; Value="18"},
    @{Cell="D5"; Value="0.07812"},
    @{Cell="E5"; Value="4.56%"},
    @{Cell="G5"; Value="18"},
    @{Cell="D6"; Value="2.254"},
    @{Cell="E6"; Value="0.27%"},
    @{Cell="G6"; Value="18"},
    @{Cell="D7"; Value="8.116"},
    @{Cell="E7"; Value="4.57%"},
    @{Cell="G7"; Value="18"},
    @{Cell="D8"; Value="4.015"},
    @{Cell="E8"; Value="6.46%"},
    @{Cell="G8"; Value="18"},
    @{Cell="D9"; Value="0.9302"},
    @{Cell="E9"; Value="1.17%"},
    @{Cell="G9"; Value="18"},
    @{Cell="D10"; Value="0.09765"},
    @{Cell="E10"; Value="4.41%"},
    @{Cell="G10"; Value="18"},
    @{Cell="D11"; Value="0.1825"},
    @{Cell="E11"; Value="5.37%"},
    @{Cell="G11"; Value="18"},
    @{Cell="D12"; Value="0.08694"},
    @{Cell="E12"; Value="3.97%"},
    @{Cell="G12"; Value="18"},
    @{Cell="D13"; Value="0.03424"},
    @{Cell="E13"; Value="4.48%"},
    @{Cell="G13"; Value="18"},
    @{Cell="D14"; Value="0.09922"},
    @{Cell="E14"; Value="-0.28%"},
    @{Cell="G14"; Value="18"},
    @{Cell="D15"; Value="0.001488"},
    @{Cell="E15"; Value="-0.44%"},
    @{Cell="G15"; Value="18"},
    @{Cell="D16"; Value="0.005659"},
    @{Cell="E16"; Value="-0.90%"},
    @{Cell="G16"; Value="18"},
    @{Cell="D17"; Value="3.494"},
    @{Cell="E17"; Value="0.58%"},
    @{Cell="G17"; Value="18"},
    @{Cell="D18"; Value="2.144"},
    @{Cell="E18"; Value="0.55%"},
    @{Cell="G18"; Value="18"},
    @{Cell="E19"; Value="2.94%"},
    @{Cell="G19"; Value="18"},
    @{Cell="D20"; Value="0.1321"},
    @{Cell="E20"; Value="1.30%"},
    @{Cell="G20"; Value="18"},
    @{Cell="D21"; Value="4.551"},
    @{Cell="E21"; Value="11.25%"},
    @{Cell="G21"; Value="18"},
    @{Cell="D22"; Value="0.2238"},
    @{Cell="E22"; Value="5.67%"},
    @{Cell="G22"; Value="18"},
    @{Cell="D23"; Value="0.04676"},
    @{Cell="E23"; Value="3.33%"},
    @{Cell="G23"; Value="18"},
    @{Cell="D24"; Value="0.001240"},
    @{Cell="E24"; Value="1.72%"},
    @{Cell="G24"; Value="18"},
    @{Cell="D25"; Value="0.004486"},
    @{Cell="E25"; Value="5.23%"},
    @{Cell="G25"; Value="18"},
    @{Cell="E26"; Value="0.47%"},
    @{Cell="G26"; Value="18"},
    @{Cell="D27"; Value="0.0002702"},
    @{Cell="E27"; Value="-20.40%"},
    @{Cell="G27"; Value="18"},
    @{Cell="G28"; Value="18"},
    @{Cell="G29"; Value="18"},
    @{Cell="G30"; Value="18"},
    @{Cell="G31"; Value="18"},
    @{Cell="G32"; Value="18"},
    @{Cell="G33"; Value="18"},
    @{Cell="G34"; Value="18"},
    @{Cell="G35"; Value="18"},
    @{Cell="G36"; Value="18"},
    @{Cell="G37"; Value="18"},
    @{Cell="G38"; Value="18"},
    @{Cell="D39"; Value="0.01753"},
    @{Cell="E39"; Value="8.36%"},
    @{Cell="G39"; Value="18"},
    @{Cell="D40"; Value="0.04710"},
    @{Cell="E40"; Value="2.98%"},
    @{Cell="G40"; Value="18"},
    @{Cell="D41"; Value="0.007845"},
    @{Cell="E41"; Value="5.34%"},
    @{Cell="G41"; Value="18"},
    @{Cell="D42"; Value="0.1413"},
    @{Cell="E42"; Value="3.97%"},
    @{Cell="G42"; Value="18"},
    @{Cell="D43"; Value="0.008741"},
    @{Cell="E43"; Value="-11.14%"},
    @{Cell="G43"; Value="18"},
    @{Cell="D44"; Value="0.002215"},
    @{Cell="E44"; Value="2.83%"},
    @{Cell="G44"; Value="18"},
    @{Cell="D45"; Value="0.009222"},
    @{Cell="E45"; Value="3.20%"},
    @{Cell="G45"; Value="18"},
    @{Cell="D46"; Value="0.00006137"},
    @{Cell="E46"; Value="0.54%"},
    @{Cell="G46"; Value="18"},
    @{Cell="E47"; Value="0.38%"},
    @{Cell="G47"; Value="18"},
    @{Cell="D48"; Value="4.073"},
    @{Cell="E48"; Value="59.69%"},
    @{Cell="G48"; Value="18"},
    @{Cell="D49"; Value="0.002693"},
    @{Cell="E49"; Value="34.65%"},
    @{Cell="G49"; Value="18"},
    @{Cell="D50"; Value="0.00002102"},
    @{Cell="E50"; Value="0.38%"},
    @{Cell="G50"; Value="18"},
    @{Cell="D51"; Value="0.0002002"},
    @{Cell="E51"; Value="0.38%"},
    @{Cell="G51"; Value="18"}
)

foreach ($edit in $edits) {
    $cell = $edit.Cell
    $row = [int]($cell -replace '[A-Za-z]+', '')
    $formatDonor = "F" + $row
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $edit.Value
    $ws.Range($formatDonor).Copy()
    $ws.Range($cell).PasteSpecial(-4122)
}
